$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-4: column F becomes the literal string "feedback" (set first so it
# lands at the same shared-string slot the original author used).
$ws.Range("F2").Value = "feedback"
$ws.Range("F3").Value = "feedback"
$ws.Range("F4").Value = "feedback"

# Row 1: turn into the real question header row (question text, 4 answer
# options, and the explanatory feedback text) - filled in the same order
# the author appears to have typed them.
$ws.Range("A1").Value = "What are the three buses in a Von Neumann CPU?"
$ws.Range("C1").Value = "Data, Memory, Control"
$ws.Range("D1").Value = "Control, Register, Memory"
$ws.Range("B1").Value = "Data,  Address,  Control"
$ws.Range("E1").Value = "Fetch, Decode, Execute"
$ws.Range("F1").Value = "The data bus retrieves data and instructions from main memory. The address bus sends addresses to main memory. The control bus sends read right signals to main memory"

# Rows 1-4: H goes 1 -> 0, J goes 2 -> 1 (G, I, K are unchanged)
$ws.Range("H1:H4").Value = 0
$ws.Range("J1:J4").Value = 1

# New column L, mirroring column K (value 2) for every data row
$ws.Range("L1:L4").Value = 2

# Update the active selection shown in the saved view to F1
$ws.Range("F1").Select()
